$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (TruckID, AssignedDockPosition, start_loading_time, end_loading_time)
# for rows 2 through 11 (data rows, row 1 is header)
$data = @(
    @(6,2,5,5),
    @(6,2,10,10),
    @(1,3,5,5),
    @(2,3,10,10),
    @(3,3,15,15),
    @(7,3,20,20),
    @(8,3,25,26),
    @(4,4,5,6),
    @(5,4,11,12),
    @(8,4,17,18)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
